# 自动更新Excel文件 - 2026-02-21 23:16:51
# Decrement column E (剩余 / "remaining") by 1 for every data row (rows 2-99),
# except row 36 which keeps its original value (its F column has a malformed
# date, so it was excluded from the refresh in the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # Column E
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current - 1
    }
}
